$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 216; this pushes the existing rows
# 216:231 down to 220:235 (identical content, just shifted).
$ws.Rows("216:219").Insert()

# Fill the 4 newly-inserted rows with the new weekly records
# (same shape as the rows that used to occupy 216:219, but with an
# updated date, variety, origin -- and, for row 217, volume).

# Row 216 - Especial
$ws.Cells.Item(216, 1).Value = 1
$ws.Cells.Item(216, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(216, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(216, 4).Value = 44826
$ws.Cells.Item(216, 5).Value = 15
$ws.Cells.Item(216, 6).Value = "Fruta"
$ws.Cells.Item(216, 7).Value = 100108
$ws.Cells.Item(216, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(216, 9).Value = 100108005
$ws.Cells.Item(216, 10).Value = "Piña"
$ws.Cells.Item(216, 11).Value = "Pucallpa"
$ws.Cells.Item(216, 12).Value = "Especial"
$ws.Cells.Item(216, 13).Value = 200
$ws.Cells.Item(216, 14).Value = 14000
$ws.Cells.Item(216, 15).Value = 15000
$ws.Cells.Item(216, 16).Value = 14500
$ws.Cells.Item(216, 17).Value = "$/caja 10 unidades"
$ws.Cells.Item(216, 18).Value = "Bolivia"
$ws.Cells.Item(216, 19).Value = 1450
$ws.Cells.Item(216, 20).Value = 10

# Row 217 - Primera
$ws.Cells.Item(217, 1).Value = 1
$ws.Cells.Item(217, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(217, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(217, 4).Value = 44826
$ws.Cells.Item(217, 5).Value = 15
$ws.Cells.Item(217, 6).Value = "Fruta"
$ws.Cells.Item(217, 7).Value = 100108
$ws.Cells.Item(217, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(217, 9).Value = 100108005
$ws.Cells.Item(217, 10).Value = "Piña"
$ws.Cells.Item(217, 11).Value = "Pucallpa"
$ws.Cells.Item(217, 12).Value = "Primera"
$ws.Cells.Item(217, 13).Value = 270
$ws.Cells.Item(217, 14).Value = 14000
$ws.Cells.Item(217, 15).Value = 15000
$ws.Cells.Item(217, 16).Value = 14500
$ws.Cells.Item(217, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(217, 18).Value = "Bolivia"
$ws.Cells.Item(217, 19).Value = 1208
$ws.Cells.Item(217, 20).Value = 12

# Row 218 - Segunda
$ws.Cells.Item(218, 1).Value = 1
$ws.Cells.Item(218, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(218, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(218, 4).Value = 44826
$ws.Cells.Item(218, 5).Value = 15
$ws.Cells.Item(218, 6).Value = "Fruta"
$ws.Cells.Item(218, 7).Value = 100108
$ws.Cells.Item(218, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(218, 9).Value = 100108005
$ws.Cells.Item(218, 10).Value = "Piña"
$ws.Cells.Item(218, 11).Value = "Pucallpa"
$ws.Cells.Item(218, 12).Value = "Segunda"
$ws.Cells.Item(218, 13).Value = 250
$ws.Cells.Item(218, 14).Value = 14000
$ws.Cells.Item(218, 15).Value = 15000
$ws.Cells.Item(218, 16).Value = 14500
$ws.Cells.Item(218, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(218, 18).Value = "Bolivia"
$ws.Cells.Item(218, 19).Value = 1036
$ws.Cells.Item(218, 20).Value = 14

# Row 219 - Tercera
$ws.Cells.Item(219, 1).Value = 1
$ws.Cells.Item(219, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(219, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(219, 4).Value = 44826
$ws.Cells.Item(219, 5).Value = 15
$ws.Cells.Item(219, 6).Value = "Fruta"
$ws.Cells.Item(219, 7).Value = 100108
$ws.Cells.Item(219, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(219, 9).Value = 100108005
$ws.Cells.Item(219, 10).Value = "Piña"
$ws.Cells.Item(219, 11).Value = "Pucallpa"
$ws.Cells.Item(219, 12).Value = "Tercera"
$ws.Cells.Item(219, 13).Value = 270
$ws.Cells.Item(219, 14).Value = 14000
$ws.Cells.Item(219, 15).Value = 15000
$ws.Cells.Item(219, 16).Value = 14500
$ws.Cells.Item(219, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(219, 18).Value = "Bolivia"
$ws.Cells.Item(219, 19).Value = 906
$ws.Cells.Item(219, 20).Value = 16
